# "criado os arquivos de trabalho" - add the new work item rows to the
# backlog sheet (Planilha2) and highlight the "backlog-principal" section
# header in bold, matching the other section headers on that sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha2")
$ws.Activate()

# "backlog-principal" (row 24) becomes a bold section header, like the
# other section headers on this sheet (car/estrada/texto/inimigo).
$ws.Range("B24").Font.Bold = $true

# Duplicate the existing "inicialização do projeto" task and add the new
# work files created for the project.
$ws.Range("B28").Value = "inicialização do projeto"
$ws.Range("C28").Value = "x"

$ws.Range("B29").Value = "index.html(front)"
$ws.Range("C29").Value = "x"

$ws.Range("B30").Value = "style.css(front)"
$ws.Range("C30").Value = "x"

$ws.Range("B31").Value = "f1_main.js"
$ws.Range("C31").Value = "x"

$ws.Range("B32").Value = "f1_main.js"
$ws.Range("C32").Value = "x"

# Update the view state to match where the author left the selection
# after typing the new rows.
$win = $excel.ActiveWindow
$win.ScrollRow = 16
$win.ScrollColumn = 1
$ws.Range("B39").Select()
